# Rename the sole worksheet from "tranter" to "Education class".
# Excel automatically rewrites every localSheetId-scoped defined name
# (_Toc133449782/4/6) so they reference the new sheet name too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Education class"

# Move the active selection from E8 to B1 (the view was also scrolled so
# row 31 is the top-left visible row, but that pure scroll-position state
# isn't part of the persisted cell selection/content).
$ws.Range("B1").Select() | Out-Null
